$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.345840692520142
$ws.Range("B1").Value = 1.452422857284546
$ws.Range("C1").Value = 3.932491302490234
$ws.Range("D1").Value = 3.263296127319336
$ws.Range("E1").Value = 1.073126673698425
